$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt7b"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.072919
$ws.Range("H2").Value = 0.218757
$ws.Range("I2").Value = 0.1016383815134179
$ws.Range("J2").Value = 0.1016383815134179
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005445666666666667
$ws.Range("N2").Value = 0.016337
$ws.Range("O2").Value = 0.1561347172047327
$ws.Range("P2").Value = 0.1561347172047327
$ws.Range("Q2").Value = 0.0003970925676666667
$ws.Range("R2").Value = 0.003573833109
$ws.Range("S2").Value = 0.01586927995474424
$ws.Range("T2").Value = 0.01586927995474424

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt7b"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.072919
$ws.Range("H3").Value = 0.218757
$ws.Range("I3").Value = 0.1016383815134179
$ws.Range("J3").Value = 0.1016383815134179
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.02943233333333333
$ws.Range("N3").Value = 0.088297
$ws.Range("O3").Value = 0.8438652827952674
$ws.Range("P3").Value = 0.8438652827952673
$ws.Range("Q3").Value = 0.002146176314333333
$ws.Range("R3").Value = 0.019315586829
$ws.Range("S3").Value = 0.0857691015586737
$ws.Range("T3").Value = 0.08576910155867369

# Row 4
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Wnt7b"
$ws.Range("C4").Value = "Fzd10"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6445166666666666
$ws.Range("H4").Value = 1.93355
$ws.Range("I4").Value = 0.8983616184865821
$ws.Range("J4").Value = 0.898361618486582
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005445666666666667
$ws.Range("N4").Value = 0.016337
$ws.Range("O4").Value = 0.1561347172047327
$ws.Range("P4").Value = 0.1561347172047327
$ws.Range("Q4").Value = 0.003509822927777778
$ws.Range("R4").Value = 0.03158840635
$ws.Range("S4").Value = 0.1402654372499885
$ws.Range("T4").Value = 0.1402654372499884

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt7b"
$ws.Range("C5").Value = "Fzd10"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6445166666666666
$ws.Range("H5").Value = 1.93355
$ws.Range("I5").Value = 0.8983616184865821
$ws.Range("J5").Value = 0.898361618486582
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.02943233333333333
$ws.Range("N5").Value = 0.088297
$ws.Range("O5").Value = 0.8438652827952674
$ws.Range("P5").Value = 0.8438652827952673
$ws.Range("Q5").Value = 0.01896962937222222
$ws.Range("R5").Value = 0.17072666435
$ws.Range("S5").Value = 0.7580961812365937
$ws.Range("T5").Value = 0.7580961812365936
